# 10.b.1.1 workbook update
# - The indicator code in the header row is bumped from "10.b.1" to
#   "10.b.1.1" for the Kyrgyz (A1) and English (C1) labels, matching the
#   Russian label (B1) which already reads "10.b.1.1".
# - The active selection is moved to L8 (an empty cell just past the used
#   range), matching the saved cursor position recorded in the file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "10.b.1.1 Агымдардын түрлөрү жана алуучу өлкөлөр жана донор-өлкөлөр боюнча бөлунүшүндөгү  өнүктүрүү максатында ресурстар агымынын жалпы көлөмү"
$ws.Range("C1").Value = "10.b.1.1 Total resource flows for development, by recipient and donor countries and type of flow (e.g. official development assistance, foreign direct investment and other flows)"

$ws.Range("L8").Select()
